$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.642.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").Value = "'1.761.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.01%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'324.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "'0.4273"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.78%  "
$ws.Range("D8").Value = "'0.3611"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.71%  "
$ws.Range("D9").Value = "'0.07594"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("D10").Value = "'42.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.68%  "
$ws.Range("D11").Value = "'1.109"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.61%  "
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "'20.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.13%  "
$ws.Range("D14").Value = "'6.077"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.99%  "
$ws.Range("D15").Value = "'7.248"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.94%  "
$ws.Range("D16").Value = "'1.764.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.07%  "
$ws.Range("D17").Value = "'92.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "'0.00001069"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("D19").Value = "'0.06443"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").Value = "'17.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.09%  "
$ws.Range("D22").Value = "'5.902"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.89%  "
$ws.Range("D23").Value = "'27.685.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.25%  "
$ws.Range("E24").Value = "  -3.14%  "
$ws.Range("D25").Value = "'2.103"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.79%  "
$ws.Range("D26").Value = "'162.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").Value = "'20.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("D28").Value = "'1.961.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.88%  "
$ws.Range("D29").Value = "'2.162"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.87%  "
$ws.Range("D30").Value = "'125.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'1.105"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.51%  "
$ws.Range("D32").Value = "'5.625"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.31%  "
$ws.Range("D33").Value = "'3.674"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.35%  "
$ws.Range("D34").Value = "'0.08949"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("D35").Value = "'12.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.63%  "
$ws.Range("D36").Value = "'0.02308"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.01%  "
$ws.Range("D37").Value = "'0.2115"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'0.06022"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.91%  "
$ws.Range("D39").Value = "'0.6360"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("D40").Value = "'4.966"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.82%  "
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("D42").Value = "'0.9999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E43").Value = "  -2.90%  "
$ws.Range("D44").Value = "'7.918"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.06%  "
$ws.Range("D45").Value = "'13.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.24%  "
$ws.Range("D46").Value = "'0.5957"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.72%  "
$ws.Range("D47").Value = "'3.706"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D49").Value = "'123.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.07%  "
$ws.Range("D50").Value = "'1.170"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("D51").Value = "'0.06877"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.91%  "
